$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price column (D) to Text format so numeric-looking strings
# (e.g. "24.542.09", "0.3260") are preserved exactly, matching the
# source data which stores prices as literal text.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '24.542.09'
$ws.Range("E2").Value = '  -0.28%  '

# Row 3
$ws.Range("D3").Value = '1.660.35'
$ws.Range("E3").Value = '  -2.48%  '

# Row 4
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  -0.12%  '

# Row 5
$ws.Range("D5").Value = '307.41'
$ws.Range("E5").Value = '  -0.11%  '

# Row 6
$ws.Range("D6").Value = '0.9963'
$ws.Range("E6").Value = '  -0.01%  '

# Row 7
$ws.Range("D7").Value = '0.3616'
$ws.Range("E7").Value = '  -2.77%  '

# Row 8
$ws.Range("D8").Value = '47.73'
$ws.Range("E8").Value = '  -2.24%  '

# Row 9
$ws.Range("D9").Value = '0.3260'
$ws.Range("E9").Value = '  -5.04%  '

# Row 10
$ws.Range("D10").Value = '1.126'
$ws.Range("E10").Value = '  -4.42%  '

# Row 11
$ws.Range("D11").Value = '0.06953'
$ws.Range("E11").Value = '  -6.28%  '

# Row 12
$ws.Range("D12").Value = '0.9958'
$ws.Range("E12").Value = '  +0.10%  '

# Row 13
$ws.Range("D13").Value = '5.908'
$ws.Range("E13").Value = '  -4.77%  '

# Row 14
$ws.Range("D14").Value = '19.40'
$ws.Range("E14").Value = '  -6.80%  '

# Row 15
$ws.Range("B15").Value = 'WrappedEther'
$ws.Range("C15").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D15").Value = '1.664.50'
$ws.Range("E15").Value = '  -2.40%  '

# Row 16
$ws.Range("B16").Value = 'Chainlink'
$ws.Range("C16").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D16").Value = '6.555'
$ws.Range("E16").Value = '  -5.08%  '

# Row 17
$ws.Range("D17").Value = '0.00001046'
$ws.Range("E17").Value = '  -6.23%  '

# Row 18
$ws.Range("D18").Value = '0.06522'
$ws.Range("E18").Value = '  -2.33%  '

# Row 19
$ws.Range("D19").Value = '0.9965'
$ws.Range("E19").Value = '  -0.02%  '

# Row 20
$ws.Range("D20").Value = '76.64'
$ws.Range("E20").Value = '  -7.73%  '

# Row 21
$ws.Range("D21").Value = '5.923'
$ws.Range("E21").Value = '  -6.28%  '

# Row 22
$ws.Range("D22").Value = '15.70'
$ws.Range("E22").Value = '  -7.71%  '

# Row 23
$ws.Range("D23").Value = '12.72'
$ws.Range("E23").Value = '  -2.90%  '

# Row 24
$ws.Range("D24").Value = '24.576.59'
$ws.Range("E24").Value = '  -0.25%  '

# Row 25
$ws.Range("D25").Value = '2.453'
$ws.Range("E25").Value = '  +2.21%  '

# Row 26
$ws.Range("D26").Value = '2.320'
$ws.Range("E26").Value = '  -15.86%  '

# Row 27
$ws.Range("D27").Value = '146.72'
$ws.Range("E27").Value = '  -1.63%  '

# Row 28
$ws.Range("D28").Value = '18.46'
$ws.Range("E28").Value = '  -7.95%  '

# Row 29
$ws.Range("D29").Value = '1.846.92'
$ws.Range("E29").Value = '  -2.30%  '

# Row 30
$ws.Range("B30").Value = 'ImmutableX'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D30").Value = '1.198'
$ws.Range("E30").Value = '  +2.55%  '

# Row 31
$ws.Range("B31").Value = 'BitcoinCash'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D31").Value = '124.20'
$ws.Range("E31").Value = '  -4.98%  '

# Row 32
$ws.Range("D32").Value = '4.051'
$ws.Range("E32").Value = '  -3.21%  '

# Row 33
$ws.Range("D33").Value = '5.632'
$ws.Range("E33").Value = '  -15.86%  '

# Row 34
$ws.Range("B34").Value = 'WEMIXTOKEN'
$ws.Range("C34").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D34").Value = '1.689'
$ws.Range("E34").Value = '  -4.33%  '

# Row 35
$ws.Range("B35").Value = 'Stellar'
$ws.Range("C35").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D35").Value = '0.08357'
$ws.Range("E35").Value = '  -4.57%  '

# Row 36
$ws.Range("D36").Value = '12.39'
$ws.Range("E36").Value = '  -8.44%  '

# Row 37
$ws.Range("D37").Value = '5.203'
$ws.Range("E37").Value = '  -5.24%  '

# Row 38
$ws.Range("D38").Value = '0.06059'
$ws.Range("E38").Value = '  -6.79%  '

# Row 39
$ws.Range("B39").Value = 'TrustWalletToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D39").Value = '1.208'
$ws.Range("E39").Value = '  -4.66%  '

# Row 40
$ws.Range("B40").Value = 'Algorand'
$ws.Range("C40").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D40").Value = '0.2056'
$ws.Range("E40").Value = '  -6.89%  '

# Row 41
$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").Value = '8.213'
$ws.Range("E41").Value = '  -7.51%  '

# Row 42
$ws.Range("B42").Value = 'VeChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D42").Value = '0.02185'
$ws.Range("E42").Value = '  -7.26%  '

# Row 43
$ws.Range("D43").Value = '0.9960'
$ws.Range("E43").Value = '  +0.03%  '

# Row 44
$ws.Range("D44").Value = '0.5905'
$ws.Range("E44").Value = '  -7.24%  '

# Row 45
$ws.Range("D45").Value = '3.739'
$ws.Range("E45").Value = '  -1.36%  '

# Row 46
$ws.Range("D46").Value = '12.69'
$ws.Range("E46").Value = '  -8.32%  '

# Row 47
$ws.Range("D47").Value = '0.5593'
$ws.Range("E47").Value = '  -7.55%  '

# Row 48
$ws.Range("D48").Value = '122.56'
$ws.Range("E48").Value = '  -4.65%  '

# Row 49
$ws.Range("D49").Value = '1.939'
$ws.Range("E49").Value = '  -7.73%  '

# Row 50
$ws.Range("D50").Value = '0.06933'
$ws.Range("E50").Value = '  -4.37%  '

# Row 51
$ws.Range("D51").Value = '74.24'
$ws.Range("E51").Value = '  -5.67%  '
